$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Copy()
$ws.Range("A2:A82").PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value = "2005Q1"
$ws.Cells.Item(3, 1).Value = "2005Q2"
$ws.Cells.Item(4, 1).Value = "2005Q3"
$ws.Cells.Item(5, 1).Value = "2005Q4"
$ws.Cells.Item(6, 1).Value = "2006Q1"
$ws.Cells.Item(7, 1).Value = "2006Q2"
$ws.Cells.Item(8, 1).Value = "2006Q3"
$ws.Cells.Item(9, 1).Value = "2006Q4"
$ws.Cells.Item(10, 1).Value = "2007Q1"
$ws.Cells.Item(11, 1).Value = "2007Q2"
$ws.Cells.Item(12, 1).Value = "2007Q3"
$ws.Cells.Item(13, 1).Value = "2007Q4"
$ws.Cells.Item(14, 1).Value = "2008Q1"
$ws.Cells.Item(15, 1).Value = "2008Q2"
$ws.Cells.Item(16, 1).Value = "2008Q3"
$ws.Cells.Item(17, 1).Value = "2008Q4"
$ws.Cells.Item(18, 1).Value = "2009Q1"
$ws.Cells.Item(19, 1).Value = "2009Q2"
$ws.Cells.Item(20, 1).Value = "2009Q3"
$ws.Cells.Item(21, 1).Value = "2009Q4"
$ws.Cells.Item(22, 1).Value = "2010Q1"
$ws.Cells.Item(23, 1).Value = "2010Q2"
$ws.Cells.Item(24, 1).Value = "2010Q3"
$ws.Cells.Item(25, 1).Value = "2010Q4"
$ws.Cells.Item(26, 1).Value = "2011Q1"
$ws.Cells.Item(27, 1).Value = "2011Q3"
$ws.Cells.Item(28, 1).Value = "2011Q4"
$ws.Cells.Item(29, 1).Value = "2012Q1"
$ws.Cells.Item(30, 1).Value = "2012Q2"
$ws.Cells.Item(31, 1).Value = "2012Q3"
$ws.Cells.Item(32, 1).Value = "2012Q4"
$ws.Cells.Item(33, 1).Value = "2013Q1"
$ws.Cells.Item(34, 1).Value = "2013Q2"
$ws.Cells.Item(35, 1).Value = "2013Q3"
$ws.Cells.Item(36, 1).Value = "2013Q4"
$ws.Cells.Item(37, 1).Value = "2014Q1"
$ws.Cells.Item(38, 1).Value = "2014Q3"
$ws.Cells.Item(39, 1).Value = "2014Q4"
$ws.Cells.Item(40, 1).Value = "2015Q1"
$ws.Cells.Item(41, 1).Value = "2015Q2"
$ws.Cells.Item(42, 1).Value = "2015Q3"
$ws.Cells.Item(43, 1).Value = "2015Q4"
$ws.Cells.Item(44, 1).Value = "2016Q1"
$ws.Cells.Item(45, 1).Value = "2016Q2"
$ws.Cells.Item(46, 1).Value = "2016Q3"
$ws.Cells.Item(47, 1).Value = "2016Q4"
$ws.Cells.Item(48, 1).Value = "2017Q1"
$ws.Cells.Item(49, 1).Value = "2017Q2"
$ws.Cells.Item(50, 1).Value = "2017Q3"
$ws.Cells.Item(51, 1).Value = "2017Q4"
$ws.Cells.Item(52, 1).Value = "2018Q1"
$ws.Cells.Item(53, 1).Value = "2018Q2"
$ws.Cells.Item(54, 1).Value = "2018Q3"
$ws.Cells.Item(55, 1).Value = "2018Q4"
$ws.Cells.Item(56, 1).Value = "2019Q1"
$ws.Cells.Item(57, 1).Value = "2019Q2"
$ws.Cells.Item(58, 1).Value = "2019Q3"
$ws.Cells.Item(59, 1).Value = "2019Q4"
$ws.Cells.Item(60, 1).Value = "2020Q1"
$ws.Cells.Item(61, 1).Value = "2020Q2"
$ws.Cells.Item(62, 1).Value = "2020Q3"
$ws.Cells.Item(63, 1).Value = "2020Q4"
$ws.Cells.Item(64, 1).Value = "2021Q1"
$ws.Cells.Item(65, 1).Value = "2021Q2"
$ws.Cells.Item(66, 1).Value = "2021Q3"
$ws.Cells.Item(67, 1).Value = "2021Q4"
$ws.Cells.Item(68, 1).Value = "2022Q1"
$ws.Cells.Item(69, 1).Value = "2022Q2"
$ws.Cells.Item(70, 1).Value = "2022Q3"
$ws.Cells.Item(71, 1).Value = "2022Q4"
$ws.Cells.Item(72, 1).Value = "2023Q1"
$ws.Cells.Item(73, 1).Value = "2023Q2"
$ws.Cells.Item(74, 1).Value = "2023Q3"
$ws.Cells.Item(75, 1).Value = "2023Q4"
$ws.Cells.Item(76, 1).Value = "2024Q1"
$ws.Cells.Item(77, 1).Value = "2024Q2"
$ws.Cells.Item(78, 1).Value = "2024Q3"
$ws.Cells.Item(79, 1).Value = "2024Q4"
$ws.Cells.Item(80, 1).Value = "2025Q1"
$ws.Cells.Item(81, 1).Value = "2025Q2"
$ws.Cells.Item(82, 1).Value = "2025Q3"
